# Refresh TPM-derived NATMI metrics (Pdcd1lg2-Rgmb) with newly computed values.
# Each entry below is a worksheet range paired with the row of replacement values
# that Excel will write in one shot (mirrors re-running the analysis with new TPM input).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Range = "G2:J2"; Values = @(0.6193099999999999, 1.85793, 0.0433418173432821, 0.0433418173432821) },
    @{ Range = "M2:T2"; Values = @(12.8066855, 25.613371, 0.3362022094228136, 0.2862725115007782, 7.931308397004999, 47.58785038203, 0.01457161475121146, 0.01240757090386935) },
    @{ Range = "G3:J3"; Values = @(0.6193099999999999, 1.85793, 0.0433418173432821, 0.0433418173432821) },
    @{ Range = "M3"; Values = @(5.556209333333332) },
    @{ Range = "O3:T3"; Values = @(0.1458620853836322, 0.186299960314954, 3.441016002226666, 30.96914402003999, 0.006321927862007606, 0.00807457885103144) },
    @{ Range = "G4:J4"; Values = @(0.6193099999999999, 1.85793, 0.0433418173432821, 0.0433418173432821) },
    @{ Range = "M4:T4"; Values = @(3.618423666666667, 10.855271, 0.0949911693670569, 0.121325915756718, 2.240925961003333, 20.16833364903, 0.004117089911931754, 0.005258485679734104) },
    @{ Range = "G5:J5"; Values = @(0.6193099999999999, 1.85793, 0.0433418173432821, 0.0433418173432821) },
    @{ Range = "M5:T5"; Values = @(11.9979585, 23.995917, 0.3149714386492294, 0.2681947419320252, 7.430455678634998, 44.58273407180999, 0.01365143456228569, 0.01162404751724652) },
    @{ Range = "G6:T6"; Values = @(0.6193099999999999, 1.85793, 0.0433418173432821, 0.0433418173432821, 2, 0.6666666666666666, 0.6145123333333333, 1.843537, 0.01613223063721173, 0.02060462744379138, 0.3805736331566666, 3.42516269841, 0.00069920019361773, 0.0008930419990951834) },
    @{ Range = "G7:J7"; Values = @(0.6193099999999999, 1.85793, 0.0433418173432821, 0.0433418173432821) },
    @{ Range = "M7:T7"; Values = @(3.498421666666667, 10.495265, 0.09184086654005638, 0.1173022430517332, 2.166607522383333, 19.49946770145, 0.003980550062227872, 0.005084092392305503) },
    @{ Range = "E8:J8"; Values = @(3, 1, 13.47892, 40.43676, 0.9433093097555539, 0.9433093097555537) },
    @{ Range = "M8:T8"; Values = @(12.8066855, 25.613371, 0.3362022094228136, 0.2862725115007782, 172.62028931966, 1035.72173591796, 0.3171426741089265, 0.2700435252257879) },
    @{ Range = "E9:J9"; Values = @(3, 1, 13.47892, 40.43676, 0.9433093097555539, 0.9433093097555537) },
    @{ Range = "M9"; Values = @(5.556209333333332) },
    @{ Range = "O9:T9"; Values = @(0.1458620853836322, 0.186299960314954, 74.89170110725333, 674.0253099652799, 0.1375930630827398, 0.1757384869721863) },
    @{ Range = "E10:J10"; Values = @(3, 1, 13.47892, 40.43676, 0.9433093097555539, 0.9433093097555537) },
    @{ Range = "M10:T10"; Values = @(3.618423666666667, 10.855271, 0.0949911693670569, 0.121325915756718, 48.77244312910668, 438.9519881619601, 0.08960605440851137, 0.1144478658479301) },
    @{ Range = "E11:J11"; Values = @(3, 1, 13.47892, 40.43676, 0.9433093097555539, 0.9433093097555537) },
    @{ Range = "M11:T11"; Values = @(11.9979585, 23.995917, 0.3149714386492294, 0.2681947419320252, 161.71952278482, 970.31713670892, 0.2971154903849184, 0.2529905968919676) },
    @{ Range = "E12:T12"; Values = @(3, 1, 13.47892, 40.43676, 0.9433093097555539, 0.9433093097555537, 2, 0.6666666666666666, 0.6145123333333333, 1.843537, 0.01613223063721173, 0.02060462744379138, 8.282962580013333, 74.54666322012, 0.01521768334720559, 0.01943653689177318) },
    @{ Range = "E13:J13"; Values = @(3, 1, 13.47892, 40.43676, 0.9433093097555539, 0.9433093097555537) },
    @{ Range = "M13:T13"; Values = @(3.498421666666667, 10.495265, 0.09184086654005638, 0.1173022430517332, 47.15494577126667, 424.3945119414, 0.08663434442325253, 0.1106522979259087) },
    @{ Range = "G14:J14"; Values = @(0.1907416666666667, 0.572225, 0.01334887290116399, 0.01334887290116398) },
    @{ Range = "M14:T14"; Values = @(12.8066855, 25.613371, 0.3362022094228136, 0.2862725115007782, 2.442768536745834, 14.656611220475, 0.004487920562675655, 0.003821415371120893) },
    @{ Range = "G15:J15"; Values = @(0.1907416666666667, 0.572225, 0.01334887290116399, 0.01334887290116398) },
    @{ Range = "M15"; Values = @(5.556209333333332) },
    @{ Range = "O15:T15"; Values = @(0.1458620853836322, 0.186299960314954, 1.059800628588889, 9.538205657299999, 0.001947094438884836, 0.002486894491736215) },
    @{ Range = "G16:J16"; Values = @(0.1907416666666667, 0.572225, 0.01334887290116399, 0.01334887290116398) },
    @{ Range = "M16:T16"; Values = @(3.618423666666667, 10.855271, 0.0949911693670569, 0.121325915756718, 0.6901841608861112, 6.211657447975001, 0.001268025046613785, 0.001619564229053758) },
    @{ Range = "G17:J17"; Values = @(0.1907416666666667, 0.572225, 0.01334887290116399, 0.01334887290116398) },
    @{ Range = "M17:T17"; Values = @(11.9979585, 23.995917, 0.3149714386492294, 0.2681947419320252, 2.2885106008875, 13.731063605325, 0.004204513702025334, 0.00358009752281108) },
    @{ Range = "G18:T18"; Values = @(0.1907416666666667, 0.572225, 0.01334887290116399, 0.01334887290116398, 2, 0.6666666666666666, 0.6145123333333333, 1.843537, 0.01613223063721173, 0.02060462744379138, 0.1172131066472222, 1.054917959825, 0.000215347096388403, 0.0002750485529230064) },
    @{ Range = "G19:J19"; Values = @(0.1907416666666667, 0.572225, 0.01334887290116399, 0.01334887290116398) },
    @{ Range = "M19:T19"; Values = @(3.498421666666667, 10.495265, 0.09184086654005638, 0.1173022430517332, 0.6672947794027778, 6.005653014625, 0.001225972054575977, 0.1144478658479301) }
)

foreach ($update in $updates) {
    $values = $update.Values
    if ($values.Count -eq 1) {
        $ws.Range($update.Range).Value = $values[0]
    } else {
        $arr = New-Object 'object[,]' 1,$values.Count
        for ($i = 0; $i -lt $values.Count; $i++) {
            $arr[0,$i] = $values[$i]
        }
        $ws.Range($update.Range).Value = $arr
    }
}
